$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,15
$data[0,0] = 0
$data[0,1] = 7.559516090224524
$data[0,2] = 5.465822424965745
$data[0,3] = 8.753894261875978
$data[0,4] = 28.0593879206004
$data[0,5] = 34.68435196145695
$data[0,6] = 2.989853608322755
$data[0,7] = 3.493842785080902
$data[0,8] = 11.92832574873758
$data[0,9] = 21.22594117044643
$data[0,10] = 6.259161153991369
$data[0,11] = 24.0865006833597
$data[0,12] = 5.923802453973462
$data[0,13] = 0
$data[0,14] = 0
$data[1,0] = 0
$data[1,1] = 7.124092137678482
$data[1,2] = 5.199047992633684
$data[1,3] = 8.40911611810783
$data[1,4] = 27.61315632222255
$data[1,5] = 34.11056643848428
$data[1,6] = 3.270007069104465
$data[1,7] = 3.739979490149482
$data[1,8] = 11.92285127996455
$data[1,9] = 20.95561702721032
$data[1,10] = 6.143040444493806
$data[1,11] = 22.55224054580047
$data[1,12] = 5.825423186599146
$data[1,13] = 0
$data[1,14] = 0
$data[2,0] = 0
$data[2,1] = 6.840445043597139
$data[2,2] = 5.029241627451165
$data[2,3] = 8.189559202080703
$data[2,4] = 27.34593764411554
$data[2,5] = 33.77089899046396
$data[2,6] = 3.447792330957014
$data[2,7] = 3.896957523046751
$data[2,8] = 11.92376780044085
$data[2,9] = 20.79283126414193
$data[2,10] = 6.069532740238312
$data[2,11] = 21.55662122484144
$data[2,12] = 5.764637577740717
$data[2,13] = 0
$data[2,14] = 0
$data[3,0] = 0
$data[3,1] = 6.712059923103203
$data[3,2] = 4.960934889908174
$data[3,3] = 8.097056164554422
$data[3,4] = 27.22434501324423
$data[3,5] = 33.61089279563299
$data[3,6] = 3.522263225539638
$data[3,7] = 3.96514814789859
$data[3,8] = 11.92066471530391
$data[3,9] = 20.71512757117475
$data[3,10] = 6.037833440821591
$data[3,11] = 21.1375484347812
$data[3,12] = 5.740911469489582
$data[3,13] = 0
$data[3,14] = 0
$data[4,0] = 0
$data[4,1] = 6.679854369616709
$data[4,2] = 4.952379379516024
$data[4,3] = 8.080219296521275
$data[4,4] = 27.18661819278086
$data[4,5] = 33.55425724822834
$data[4,6] = 3.535380122094998
$data[4,7] = 3.980039277161069
$data[4,8] = 11.91469967275135
$data[4,9] = 20.6873866182754
$data[4,10] = 6.031074469698384
$data[4,11] = 21.06725959852739
$data[4,12] = 5.738355358134902
$data[4,13] = 0
$data[4,14] = 0
$data[5,0] = 0
$data[5,1] = 6.810237593976755
$data[5,2] = 5.036029653949347
$data[5,3] = 8.184648599543639
$data[5,4] = 27.29636755264428
$data[5,5] = 33.68651480876624
$data[5,6] = 3.45057173526129
$data[5,7] = 3.90712717288565
$data[5,8] = 11.90869473601057
$data[5,9] = 20.75128080980885
$data[5,10] = 6.065146019174052
$data[5,11] = 21.55128867351768
$data[5,12] = 5.768090442279433
$data[5,13] = 0
$data[5,14] = 0
$data[6,0] = 0
$data[6,1] = 7.377539156045561
$data[6,2] = 5.384689568192986
$data[6,3] = 8.632049394688904
$data[6,4] = 27.84203271611286
$data[6,5] = 34.37731502034536
$data[6,6] = 3.087525601869366
$data[6,7] = 3.589474354855306
$data[6,8] = 11.90577002319867
$data[6,9] = 21.07933998646085
$data[6,10] = 6.214492746635216
$data[6,11] = 23.56884687415242
$data[6,12] = 5.894874390584333
$data[6,13] = 0
$data[6,14] = 0
$data[7,0] = 0
$data[7,1] = 8.407130243170013
$data[7,2] = 6.006333429998914
$data[7,3] = 9.451041283385408
$data[7,4] = 29.04802701405723
$data[7,5] = 35.97614514725242
$data[7,6] = 2.416979769289911
$data[7,7] = 2.996610419393182
$data[7,8] = 11.96267368436469
$data[7,9] = 21.81866262352322
$data[7,10] = 6.49637339269645
$data[7,11] = 27.09967842119948
$data[7,12] = 6.133571970635619
$data[7,13] = 0
$data[7,14] = 0
$data[8,0] = 0
$data[8,1] = 9.036457782735646
$data[8,2] = 6.416232312642163
$data[8,3] = 9.909980430430537
$data[8,4] = 29.70782819974371
$data[8,5] = 36.81890131814104
$data[8,6] = 1.982087060352266
$data[8,7] = 2.602256348366554
$data[8,8] = 11.9528938657101
$data[8,9] = 22.19770413562594
$data[8,10] = 6.643357479822654
$data[8,11] = 29.43540262777505
$data[8,12] = 6.25491847021225
$data[8,13] = 0
$data[8,14] = 0
$data[9,0] = 0
$data[9,1] = 8.95100468452225
$data[9,2] = 6.44889577607326
$data[9,3] = 9.269308163466249
$data[9,4] = 27.92545799812547
$data[9,5] = 34.03046544920772
$data[9,6] = 2.88229931354139
$data[9,7] = 2.515349121096044
$data[9,8] = 11.3382321263616
$data[9,9] = 20.93864552720473
$data[9,10] = 6.327009440076889
$data[9,11] = 30.43982046069266
$data[9,12] = 5.881277681019414
$data[9,13] = 0
$data[9,14] = 0
$data[10,0] = 0
$data[10,1] = 8.783963890671529
$data[10,2] = 6.378736524468717
$data[10,3] = 8.664712378169437
$data[10,4] = 26.33784190816912
$data[10,5] = 31.59174531303732
$data[10,6] = 4.228118512557674
$data[10,7] = 2.504825043611858
$data[10,8] = 10.8472781617967
$data[10,9] = 19.85375223506002
$data[10,10] = 6.07982598708607
$data[10,11] = 30.81079089712379
$data[10,12] = 5.577025974284325
$data[10,13] = 0
$data[10,14] = 0
$data[11,0] = 0
$data[11,1] = 8.504211366807986
$data[11,2] = 6.232733028234057
$data[11,3] = 8.042522007657483
$data[11,4] = 24.7228834032151
$data[11,5] = 29.13307846370278
$data[11,6] = 5.645030539953513
$data[11,7] = 2.554205295955391
$data[11,8] = 10.4012262176772
$data[11,9] = 18.77411873035082
$data[11,10] = 5.868964535554615
$data[11,11] = 30.72956680685675
$data[11,12] = 5.319123353472177
$data[11,13] = 0
$data[11,14] = 0
$data[12,0] = 0
$data[12,1] = 8.247697150647765
$data[12,2] = 6.096009333704134
$data[12,3] = 7.60539872301466
$data[12,4] = 23.57356785640574
$data[12,5] = 27.3907899691164
$data[12,6] = 6.655728352477559
$data[12,7] = 2.627788425910068
$data[12,8] = 10.11132105563777
$data[12,9] = 18.01735678000222
$data[12,10] = 5.747017772939114
$data[12,11] = 30.467308687607
$data[12,12] = 5.171502122738233
$data[12,13] = 0
$data[12,14] = 0
$data[13,0] = 0
$data[13,1] = 8.151761311287183
$data[13,2] = 6.048915757143341
$data[13,3] = 7.490391651664166
$data[13,4] = 23.27568634251819
$data[13,5] = 26.94170884770359
$data[13,6] = 6.893629082297322
$data[13,7] = 2.666502354336215
$data[13,8] = 10.04627613792728
$data[13,9] = 17.82297651017181
$data[13,10] = 5.71844864391392
$data[13,11] = 30.30598466501647
$data[13,12] = 5.140255171780748
$data[13,13] = 0
$data[13,14] = 0
$data[14,0] = 0
$data[14,1] = 7.913248063358745
$data[14,2] = 5.904021994130336
$data[14,3] = 7.394379409231984
$data[14,4] = 23.26581861146661
$data[14,5] = 26.98095610417067
$data[14,6] = 6.721861002399367
$data[14,7] = 2.829243346657083
$data[14,8] = 10.13615963667838
$data[14,9] = 17.84688042458122
$data[14,10] = 5.693431045872221
$data[14,11] = 29.36651319044827
$data[14,12] = 5.139404006568423
$data[14,13] = 0
$data[14,14] = 0
$data[15,0] = 0
$data[15,1] = 7.868849088841869
$data[15,2] = 5.86742048269936
$data[15,3] = 7.552186771050144
$data[15,4] = 23.86988469750043
$data[15,5] = 27.93808149284184
$data[15,6] = 6.029802506380736
$data[15,7] = 2.916061279621335
$data[15,8] = 10.35541184395499
$data[15,9] = 18.2643875011064
$data[15,10] = 5.733459216886127
$data[15,11] = 28.7758366037179
$data[15,12] = 5.21180737682084
$data[15,13] = 0
$data[15,14] = 0
$data[16,0] = 0
$data[16,1] = 8.009541365909561
$data[16,2] = 5.914654864837416
$data[16,3] = 7.977824352853339
$data[16,4] = 25.11253310094139
$data[16,5] = 29.86230265372517
$data[16,6] = 4.846571668736859
$data[16,7] = 2.93326695448865
$data[16,8] = 10.72788413117864
$data[16,9] = 19.10673242253484
$data[16,10] = 5.865900648718116
$data[16,11] = 28.43043363888723
$data[16,12] = 5.385270510009907
$data[16,13] = 0
$data[16,14] = 0
$data[17,0] = 0
$data[17,1] = 8.246053961731006
$data[17,2] = 6.034935801183501
$data[17,3] = 8.614536535926302
$data[17,4] = 26.71736038299147
$data[17,5] = 32.31543465920971
$data[17,6] = 3.477809482982827
$data[17,7] = 2.910150351141169
$data[17,8] = 11.18564196149491
$data[17,9] = 20.18301552208176
$data[17,10] = 6.100012554465422
$data[17,11] = 28.31213860479293
$data[17,12] = 5.67228721854624
$data[17,13] = 0
$data[17,14] = 0
$data[18,0] = 0
$data[18,1] = 8.803746326366023
$data[18,2] = 6.329331273404806
$data[18,3] = 9.777559946342574
$data[18,4] = 29.38442957212384
$data[18,5] = 36.34179200460635
$data[18,6] = 2.098936628887856
$data[18,7] = 2.738452525640375
$data[18,8] = 11.90521603188889
$data[18,9] = 21.97189350318592
$data[18,10] = 6.592232821605092
$data[18,11] = 28.83944391848903
$data[18,12] = 6.231704805269311
$data[18,13] = 0
$data[18,14] = 0
$data[19,0] = 0
$data[19,1] = 9.339977076625011
$data[19,2] = 6.656311326255911
$data[19,3] = 10.27368499497078
$data[19,4] = 30.29633823714362
$data[19,5] = 37.61739549178271
$data[19,6] = 1.7237893302097
$data[19,7] = 2.588509586579113
$data[19,8] = 12.02343062150563
$data[19,9] = 22.54299336864825
$data[19,10] = 6.777240810983696
$data[19,11] = 30.55069329532636
$data[19,12] = 6.404377517360755
$data[19,13] = 0
$data[19,14] = 0
$data[20,0] = 0
$data[20,1] = 9.684036389461799
$data[20,2] = 6.847229602746106
$data[20,3] = 10.54245782236478
$data[20,4] = 30.83188522532251
$data[20,5] = 38.37283310762318
$data[20,6] = 1.623966134413531
$data[20,7] = 2.775746758694647
$data[20,8] = 12.08955604470509
$data[20,9] = 22.88441202815501
$data[20,10] = 6.876413340020032
$data[20,11] = 31.62300744960155
$data[20,12] = 6.484791138160535
$data[20,13] = 0
$data[20,14] = 0
$data[21,0] = 0
$data[21,1] = 9.526863626176834
$data[21,2] = 6.738928109878382
$data[21,3] = 10.40326051588798
$data[21,4] = 30.59613461526988
$data[21,5] = 38.05481234915551
$data[21,6] = 1.615269826684802
$data[21,7] = 2.671665381408537
$data[21,8] = 12.07083743018507
$data[21,9] = 22.7456620234044
$data[21,10] = 6.827787642576939
$data[21,11] = 31.05497279597036
$data[21,12] = 6.437803413135222
$data[21,13] = 0
$data[21,14] = 0
$data[22,0] = 0
$data[22,1] = 8.866014182108659
$data[22,2] = 6.32705513084358
$data[22,3] = 9.85484253375159
$data[22,4] = 29.63621131811978
$data[22,5] = 36.74157912953248
$data[22,6] = 2.078960773650089
$data[22,7] = 2.716184232422946
$data[22,8] = 11.98453872713356
$data[22,9] = 22.15847448896125
$data[22,10] = 6.633188058836868
$data[22,11] = 28.81075923005064
$data[22,12] = 6.265120152279298
$data[22,13] = 0
$data[22,14] = 0
$data[23,0] = 0
$data[23,1] = 8.098415237976486
$data[23,2] = 5.855714477382851
$data[23,3] = 9.230714181403695
$data[23,4] = 28.63420683716902
$data[23,5] = 35.39201407724322
$data[23,6] = 2.596747878594738
$data[23,7] = 3.169881877144196
$data[23,8] = 11.91643672982863
$data[23,9] = 21.54582290126555
$data[23,10] = 6.415367300590525
$data[23,11] = 30.30598466501647
$data[23,12] = 5.140255171780748
$data[23,13] = 0
$data[23,14] = 0

$ws.Range("B2:P25").Value = $data
